# fix subset method definition, remove unneeded plotSpec() override
#
# Adds a new "components" worksheet (mirroring the existing fGroups/
# formulas/compounds status sheets) after "compounds", fills it with the
# subset-method support matrix, and makes it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- add the new "components" sheet right after "compounds" -----------
$compounds = $wb.Worksheets.Item("compounds")
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $compounds)
$newSheet.Name = "components"

# --- header row ---------------------------------------------------------
$newSheet.Range("B1").Value = "as-is"
$newSheet.Range("C1").Value = "almost as-is"
$newSheet.Range("D1").Value = "implement"
$newSheet.Range("E1").Value = "not supported"
$newSheet.Range("F1").Value = "ionize"
$newSheet.Range("G1").Value = "done"

# --- data rows -----------------------------------------------------------
$newSheet.Range("A2").Value = "`$"
$newSheet.Range("B2").Value = "X"
$newSheet.Range("G2").Value = "X"

$newSheet.Range("A3").Value = "["
$newSheet.Range("C3").Value = "X"
$newSheet.Range("G3").Value = "X"

$newSheet.Range("A4").Value = "[["
$newSheet.Range("B4").Value = "X"
$newSheet.Range("G4").Value = "X"

$newSheet.Range("A5").Value = "as.data.table"
$newSheet.Range("B5").Value = "X"
$newSheet.Range("G5").Value = "X"

$newSheet.Range("A6").Value = "componentInfo"
$newSheet.Range("B6").Value = "X"
$newSheet.Range("G6").Value = "X"

$newSheet.Range("A7").Value = "componentTable"
$newSheet.Range("B7").Value = "X"
$newSheet.Range("G7").Value = "X"

$newSheet.Range("A8").Value = "consensus"
$newSheet.Range("E8").Value = "X"

$newSheet.Range("A9").Value = "filter"
$newSheet.Range("C9").Value = "X"
$newSheet.Range("G9").Value = "X"

$newSheet.Range("A10").Value = "findFGroup"
$newSheet.Range("B10").Value = "X"
$newSheet.Range("G10").Value = "X"

$newSheet.Range("A11").Value = "groupNames"
$newSheet.Range("B11").Value = "X"
$newSheet.Range("G11").Value = "X"

$newSheet.Range("A12").Value = "initialize"
$newSheet.Range("C12").Value = "X"
$newSheet.Range("G12").Value = "X"

$newSheet.Range("A13").Value = "length"
$newSheet.Range("B13").Value = "X"
$newSheet.Range("G13").Value = "X"

$newSheet.Range("A14").Value = "names"
$newSheet.Range("B14").Value = "X"
$newSheet.Range("G14").Value = "X"

$newSheet.Range("A15").Value = "plotEIC"
$newSheet.Range("B15").Value = "X"
$newSheet.Range("D15").Value = "X"
$newSheet.Range("G15").Value = "X"
$newSheet.Range("H15").Value = "Seems enough, assuming we're not planning to merge components"

$newSheet.Range("A16").Value = "plotEICHash"
$newSheet.Range("B16").Value = "X"
$newSheet.Range("G16").Value = "X"

$newSheet.Range("A17").Value = "plotSpec"
$newSheet.Range("B17").Value = "X"
$newSheet.Range("D17").Value = "X"
$newSheet.Range("G17").Value = "X"
$newSheet.Range("H17").Value = "Seems enough, assuming we're not planning to merge components"

$newSheet.Range("A18").Value = "plotSpecHash"
$newSheet.Range("B18").Value = "X"
$newSheet.Range("G18").Value = "X"

$newSheet.Range("A19").Value = "show"
$newSheet.Range("C19").Value = "X"
$newSheet.Range("G19").Value = "X"

# --- column A width, matching the other status sheets -------------------
$newSheet.Columns.Item(1).ColumnWidth = 16.140625

# --- "compounds" keeps its own remembered selection (B1:G1), but is no
#     longer the active tab ---------------------------------------------
$compounds.Select()
$compounds.Range("B1:G1").Select()

# --- the new "components" sheet becomes the active/selected tab, with
#     the cursor parked on H17 -------------------------------------------
$newSheet.Select()
$newSheet.Range("H17").Select()
